$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 65077
$ws.Range("B2").Value = "João Guilherme Pimenta"
$ws.Range("C2").Value = "Engenharia"
$ws.Range("D2").Value = "Viagem de negocios"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45085
$ws.Range("G2").Value = 6138.86

# Row 3
$ws.Range("A3").Value = 25935
$ws.Range("B3").Value = "Luísa Ramos"
$ws.Range("C3").Value = "Vendas"
$ws.Range("D3").Value = "Problemas pessoais"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 45097
$ws.Range("G3").Value = 9427.290000000001

# Row 4
$ws.Range("A4").Value = 12387
$ws.Range("B4").Value = "Léo Câmara"
$ws.Range("C4").Value = "Atendimento ao Cliente"
$ws.Range("D4").Value = "Outros"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45087
$ws.Range("G4").Value = 4829.39

# Row 5
$ws.Range("A5").Value = 65992
$ws.Range("B5").Value = "Helena Viana"
$ws.Range("C5").Value = "Marketing"
$ws.Range("D5").Value = "Problemas pessoais"
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 45082
$ws.Range("G5").Value = 5330.32

# Row 6
$ws.Range("A6").Value = 87403
$ws.Range("B6").Value = "Maria Luísa Santos"
$ws.Range("C6").Value = "P&D"
$ws.Range("D6").Value = "Viagem de negocios"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 45087
$ws.Range("G6").Value = 5618.72

# Row 7
$ws.Range("A7").Value = 40245
$ws.Range("B7").Value = "Luna Aragão"
$ws.Range("C7").Value = "Juridico"
$ws.Range("D7").Value = "Viagem de negocios"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 45103
$ws.Range("G7").Value = 2000.45

# Row 8
$ws.Range("A8").Value = 83378
$ws.Range("B8").Value = "Murilo Castro"
$ws.Range("C8").Value = "P&D"
$ws.Range("D8").Value = "Viagem de negocios"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 45096
$ws.Range("G8").Value = 3790.44

# Row 9
$ws.Range("A9").Value = 31454
$ws.Range("B9").Value = "Dr. Apollo Mendonça"
$ws.Range("C9").Value = "Vendas"
$ws.Range("D9").Value = "Doenca"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 45088
$ws.Range("G9").Value = 3472.32

# Row 10
$ws.Range("A10").Value = 32940
$ws.Range("B10").Value = "Sra. Alexia Leão"
$ws.Range("C10").Value = "Juridico"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 7
$ws.Range("F10").Value = 45093
$ws.Range("G10").Value = 4426.71

# Row 11
$ws.Range("A11").Value = 31947
$ws.Range("B11").Value = "Sr. Renan da Cruz"
$ws.Range("C11").Value = "Recursos Humanos"
$ws.Range("D11").Value = "Problemas pessoais"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 45092
$ws.Range("G11").Value = 3014.06
